$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "31.126.24"
$ws.Range("E2").Value2 = "  +1.80%  "
$ws.Range("D3").Value2 = "1.959.89"
$ws.Range("E3").Value2 = "  +2.18%  "
$ws.Range("E4").Value2 = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "247.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.76%  "
$ws.Range("E6").Value2 = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2977"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +2.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06841"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "19.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "106.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -5.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.07757"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +2.49%  "
$ws.Range("D13").Value2 = "1.936.49"
$ws.Range("E13").Value2 = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.427"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.7129"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +6.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "285.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -3.37%  "
$ws.Range("D17").Value2 = "31.133.76"
$ws.Range("E17").Value2 = "  +1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.000007779"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +1.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "5.590"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +1.49%  "
$ws.Range("E21").Value2 = "  +0.13%  "
$ws.Range("D22").Value2 = "2.193.67"
$ws.Range("E22").Value2 = "  +1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "6.602"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "9.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +5.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "168.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "20.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "2.201"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +5.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.1062"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.445"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.817"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +18.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.507"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +9.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.05029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.7704"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.166"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.02054"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.734"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "2.133"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +5.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "6.420"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +9.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.8865"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +2.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "109.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "73.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +6.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.4469"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +0.87%  "
$ws.Range("E45").Value2 = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "998.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +18.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "7.505"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +3.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.1270"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "9.411"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "36.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.2580"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +2.99%  "
